# Update "想去人数" (F column) figures for the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 14331
$wsExhibit.Range("F3").Value = 337
$wsExhibit.Range("F4").Value = 691
$wsExhibit.Range("F6").Value = 567
$wsExhibit.Range("F7").Value = 1504
$wsExhibit.Range("F8").Value = 144

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 14331
$wsAll.Range("F3").Value = 337
$wsAll.Range("F4").Value = 691
$wsAll.Range("F8").Value = 567
$wsAll.Range("F9").Value = 1504
$wsAll.Range("F11").Value = 144
